$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dat-hang-thanh-toan")

# Update the email values for rows 2, 3 and 5 (column B)
$ws.Range("B2").Value = "hahaaaa27@gmail.com"
$ws.Range("B3").Value = "hahaaa44@gmail.com"
$ws.Range("B5").Value = "hahaaa5555@gmail.com"

# Move the active selection to B5
$ws.Range("B5").Select()
